# Add an "E-Mail" column to the worksheet with mailto: hyperlinks so that
# certificates can be sent by mail (see commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Header for the new column F -------------------------------------------
$ws.Range("F1").Value = "E-Mail"
# Give the new header the same bold formatting used by the other header cells.
$ws.Range("F1").Font.Bold = $true

# -- Data for the new column -------------------------------------------------
$email = "burak.onat@tech-academy.io"
$ws.Range("F2").Value = $email
$ws.Range("F3").Value = $email
$ws.Range("F4").Value = $email
$ws.Range("F5").Value = $email

# -- Hyperlinks ---------------------------------------------------------------
# First row gets its own hyperlink (no explicit display text - the cell's
# own value is used as the link text).
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:$email")

# The remaining rows share a single hyperlink definition over the F3:F5
# range, with an explicit display text equal to the e-mail address.
$ws.Hyperlinks.Add($ws.Range("F3:F5"), "mailto:$email", "", "", $email)

# Hyperlinks.Add only stamps the built-in "Hyperlink" style onto the first
# cell of a multi-cell range; apply it to the remaining cells too so every
# e-mail cell is rendered consistently (underlined / hyperlink colour).
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"

# -- Selection ------------------------------------------------------------
# Mirrors the author's final cursor position after making the edit.
[void]$ws.Range("F19").Select()

Write-Output "E-Mail column with hyperlinks added"
